$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 452).
for ($row = 2; $row -le 452; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45179) {
        $cell.Value = 45180
    }
}
